$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Match the header styling used by the rest of row 1 (e.g. AC1) before
# writing the new header labels, so AD1:AF1 look like the other headers
# (bold, centered, bordered).
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# Add new header cells for the team record columns (AD, AE, AF)
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Fill in the team's win/loss/tie record for every data row
for ($r = 2; $r -le 50; $r++) {
    $ws.Cells.Item($r, 30).Value = 90
    $ws.Cells.Item($r, 31).Value = 72
    $ws.Cells.Item($r, 32).Value = 0
}
